# "Interview Import Me Header.csv" data was imported into rows 85-91:
# column A = Date of Check (serial date values), column C = Check Number (text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates  = @(42992, 42999, 42990, 42992, 42948, 43034, 42964)
$checks = @("81684", "81806", "81314", "81234", "81305", "81857", "81125")

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 85 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $checks[$i]
}

# Update the view so the newly imported rows are in focus.
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range("D93").Select()
